$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - swin_t
$ws.Range("A2").Value = "swin_t"
$ws.Range("F2").Value = 28.3
$ws.Range("G2").Value = 421.4442393779755

# Row 3 - convnext_tiny
$ws.Range("A3").Value = "convnext_tiny"
$ws.Range("C3").Value = 0.9952351171863368
$ws.Range("D3").Value = 0.9952351171863367
$ws.Range("F3").Value = 28.6
$ws.Range("G3").Value = 384.7073669433594

# Row 4 - efficientnet_b0
$ws.Range("A4").Value = "efficientnet_b0"
$ws.Range("C4").Value = 0.995236606212216
$ws.Range("D4").Value = 0.995236606212216
$ws.Range("F4").Value = 5.3
$ws.Range("G4").Value = 215.5847570896149

# Row 5 - efficientnet_b3
$ws.Range("A5").Value = "efficientnet_b3"
$ws.Range("C5").Value = 0.9952351171863366
$ws.Range("D5").Value = 0.9952351171863367
$ws.Range("F5").Value = 12.2
$ws.Range("G5").Value = 312.3280458450317

# Row 6 - resnet50 (name unchanged)
$ws.Range("B6").Value = 0.9952380952380953
$ws.Range("C6").Value = 0.9952351171863368
$ws.Range("D6").Value = 0.9952351171863367
$ws.Range("E6").Value = 0.995
$ws.Range("F6").Value = 25.6
$ws.Range("G6").Value = 401.7651484012604

# Row 7 - resnet101
$ws.Range("A7").Value = "resnet101"
$ws.Range("C7").Value = 0.9928526757795052
$ws.Range("D7").Value = 0.992852675779505
$ws.Range("F7").Value = 44.5
$ws.Range("G7").Value = 374.3951771259308

# Row 8 - densenet121
$ws.Range("A8").Value = "densenet121"
$ws.Range("C8").Value = 0.9904612778260312
$ws.Range("D8").Value = 0.9904612778260311
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = 211.3503749370575

# Row 9 - vit_b_16
$ws.Range("A9").Value = "vit_b_16"
$ws.Range("B9").Value = 0.9809523809523809
$ws.Range("C9").Value = 0.9809002756089702
$ws.Range("D9").Value = 0.9809002756089702
$ws.Range("E9").Value = 0.98
$ws.Range("F9").Value = 86.59999999999999
$ws.Range("G9").Value = 282.4712982177734
